# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 00:35"

# Row 4 - Estados Unidos: refreshed totals
$ws.Range("B4").Value = 1546128
$ws.Range("C4").Value = 18464
$ws.Range("D4").Value = 353648
$ws.Range("E4").Value = 1100717
$ws.Range("G4").Value = 785
$ws.Range("H4").Value = 91763

# Rows 71-72 - Camerun overtakes Azerbaiyan in ranking
$ws.Range("A71").Value = "Camerun"
$ws.Range("B71").Value = 3529
$ws.Range("C71").Value = 424
$ws.Range("D71").Value = 1567
$ws.Range("E71").Value = 1822
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 140

$ws.Range("A72").Value = "Azerbaiyan"
$ws.Range("B72").Value = 3387
$ws.Range("C72").Value = 113
$ws.Range("D72").Value = 2055
$ws.Range("E72").Value = 1292
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 40

# Row 95 - Eslovenia: recovered / critical counts updated
$ws.Range("D95").Value = 1335
$ws.Range("E95").Value = 27

# Rows 112-118 - Uruguay jumps ahead of Mali..Principado de Andorra
$ws.Range("A112").Value = "Uruguay"
$ws.Range("B112").Value = 884
$ws.Range("C112").Value = 150
$ws.Range("D112").Value = 564
$ws.Range("E112").Value = 300
$ws.Range("H112").Value = 20

$ws.Range("A113").Value = "Mali"
$ws.Range("B113").Value = 874
$ws.Range("C113").Value = 14
$ws.Range("D113").Value = 512
$ws.Range("E113").Value = 310
$ws.Range("H113").Value = 52

$ws.Range("A114").Value = "Costa Rica"
$ws.Range("B114").Value = 866
$ws.Range("C114").Value = 3
$ws.Range("D114").Value = 575
$ws.Range("E114").Value = 281
$ws.Range("H114").Value = 10

$ws.Range("A115").Value = "Burkina Faso"
$ws.Range("B115").Value = 796
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 652
$ws.Range("E115").Value = 93
$ws.Range("H115").Value = 51

$ws.Range("A116").Value = "Paraguay"
$ws.Range("B116").Value = 788
$ws.Range("C116").Value = 2
$ws.Range("D116").Value = 219
$ws.Range("E116").Value = 558
$ws.Range("H116").Value = 11

$ws.Range("A117").Value = "Zambia"
$ws.Range("B117").Value = 761
$ws.Range("C117").Value = 8
$ws.Range("D117").Value = 192
$ws.Range("E117").Value = 562
$ws.Range("H117").Value = 7

$ws.Range("A118").Value = "Principado de Andorra"
$ws.Range("B118").Value = 761
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 624
$ws.Range("E118").Value = 86
$ws.Range("H118").Value = 51
